# "Update hosts, add IgY"
#
# 1. Terminology sheet: replace the host list (A2:A6) with the new set of
#    species and add the missing "IgY" isotype at B16.
# 2. Antibodies sheet: the Host / Isotype dropdown validations need to keep
#    pointing at the (now larger) Terminology ranges.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Terminology sheet - update host species list & add IgY isotype
# ---------------------------------------------------------------------
$terminology = $wb.Worksheets.Item("Terminology")

# The sheet ships protected; unprotect so the cells can be edited, then
# restore protection afterwards.
$terminology.Unprotect()

$terminology.Range("A2").Value = "chicken (Gallus gallus)"
$terminology.Range("A3").Value = "human (Homo sapiens)"
$terminology.Range("A4").Value = "llama (Lama glama)"
$terminology.Range("A5").Value = "mouse (Mus musculus)"
$terminology.Range("A6").Value = "alpaca (Vicugna pacos)"

$terminology.Range("B16").Value = "IgY"

$terminology.Protect()

# ---------------------------------------------------------------------
# Antibodies sheet - extend the Host/Isotype dropdown source ranges to
# cover the rows added to Terminology (Host now has 5 entries, Isotype 15).
# ---------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item("Antibodies")

$hostValidation = $antibodies.Range("B2:B100").Validation
$hostValidation.Formula1 = "=Terminology!`$A`$2:`$A`$6"

$isotypeValidation = $antibodies.Range("C2:C100").Validation
$isotypeValidation.Formula1 = "=Terminology!`$B`$2:`$B`$16"
